$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update allocation results (AffectedPop) for existing rows
$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 10

# Add new row 6 with the original "Poblacion" data (report UI addition)
$ws.Range("A6").Value = $true
$ws.Range("B6").Value = "Poblacion"
$ws.Range("C6").Value = 14.9157
$ws.Range("D6").Value = 120.7672
$ws.Range("E6").Value = 1785
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = ""

# Update row 5 with the new "Calizon" allocation result
$ws.Range("B5").Value = "Calizon"
$ws.Range("C5").Value = 14.9125
$ws.Range("D5").Value = 120.753
$ws.Range("E5").Value = 2221
